# Elec sector calibration edits
# - Add a new "US values" sheet (between "Cal" and "BHRbEF") with hardcoded
#   US EPS calibration data.
# - Point a few BHRbEF rows (natural gas peaker, crude oil, municipal solid
#   waste) at the new "US values" sheet instead of hardcoded zeros.
# - Note on the About sheet that US values are used where Korea data is
#   unavailable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "US values" worksheet right before "BHRbEF" (i.e.
#    right after "Cal") and populate it with the US EPS data table.
# ---------------------------------------------------------------------
$calSheet = $wb.Worksheets.Item("Cal")
$usValues = $wb.Worksheets.Add($null, $calSheet)
$usValues.Name = "US values"

$usValues.Range("A1").Value = "Heat Rate by Electricity Fuel BTU/(MW*hour)"
$usValues.Range("B1").Value = "preexisting"
$usValues.Range("C1").Value = "preexisting nonretiring (not used in U.S. dataset)"
$usValues.Range("D1").Value = "newly built"

$usValues.Range("A2").Value = "hard coal"
$usValues.Range("B2").Value = 10511434
$usValues.Range("C2").Value = 0
$usValues.Range("D2").Value = 10375000

$usValues.Range("A3").Value = "natural gas nonpeaker"
$usValues.Range("B3").Value = 10488088
$usValues.Range("C3").Value = 7159317
$usValues.Range("D3").Value = 6516500

$usValues.Range("A4").Value = "nuclear"
$usValues.Range("B4").Value = 10442000
$usValues.Range("C4").Value = 0
$usValues.Range("D4").Value = 10442000

$usValues.Range("A5").Value = "hydro"
$usValues.Range("B5").Value = 0
$usValues.Range("C5").Value = 0
$usValues.Range("D5").Value = 0

$usValues.Range("A6").Value = "wind"
$usValues.Range("B6").Value = 0
$usValues.Range("C6").Value = 0
$usValues.Range("D6").Value = 0

$usValues.Range("A7").Value = "solar PV"
$usValues.Range("B7").Value = 0
$usValues.Range("C7").Value = 0
$usValues.Range("D7").Value = 0

$usValues.Range("A8").Value = "solar thermal"
$usValues.Range("B8").Value = 0
$usValues.Range("C8").Value = 0
$usValues.Range("D8").Value = 0

$usValues.Range("A9").Value = "biomass"
$usValues.Range("B9").Value = 9482232
$usValues.Range("C9").Value = 0
$usValues.Range("D9").Value = 9482232

$usValues.Range("A10").Value = "geothermal"
$usValues.Range("B10").Value = 0
$usValues.Range("C10").Value = 0
$usValues.Range("D10").Value = 0

$usValues.Range("A11").Value = "petroleum"
$usValues.Range("B11").Value = 10988489
$usValues.Range("C11").Value = 0
$usValues.Range("D11").Value = 10000000

$usValues.Range("A12").Value = "natural gas peaker"
$usValues.Range("B12").Value = 9381404
$usValues.Range("C12").Value = 0
$usValues.Range("D12").Value = 8902000

$usValues.Range("A13").Value = "lignite"
$usValues.Range("B13").Value = 11488776
$usValues.Range("C13").Value = 0
$usValues.Range("D13").Value = 11339657

$usValues.Range("A14").Value = "offshore wind"
$usValues.Range("B14").Value = 0
$usValues.Range("C14").Value = 0
$usValues.Range("D14").Value = 0

$usValues.Range("A15").Value = "crude oil"
$usValues.Range("B15").Value = 7713158
$usValues.Range("C15").Value = 0
$usValues.Range("D15").Value = 7713158

$usValues.Range("A16").Value = "heavy or residual fuel oil"
$usValues.Range("B16").Value = 10719153
$usValues.Range("C16").Value = 0
$usValues.Range("D16").Value = 10719153

$usValues.Range("A17").Value = "municipal solid waste"
$usValues.Range("B17").Value = 18894208
$usValues.Range("C17").Value = 0
$usValues.Range("D17").Value = 9482232

# ---------------------------------------------------------------------
# 2) Update BHRbEF to pull the rows with no Korea-specific data from the
#    new "US values" sheet instead of the old hardcoded zeros.
# ---------------------------------------------------------------------
$bhrbef = $wb.Worksheets.Item("BHRbEF")

# row 12: natural gas peaker
$bhrbef.Range("B12").Formula = "='US values'!B12"
$bhrbef.Range("C12").Formula = "='US values'!C12"
$bhrbef.Range("D12").Formula = "='US values'!D12"
# D12 previously had no explicit number format (it was a bare 0); match
# the rest of the row now that it carries a real figure.
$bhrbef.Range("D12").NumberFormat = $bhrbef.Range("C12").NumberFormat

# row 15: crude oil
$bhrbef.Range("B15").Formula = "='US values'!B15"
$bhrbef.Range("C15").Formula = "='US values'!C15"
$bhrbef.Range("D15").Formula = "='US values'!D15"

# row 17: municipal solid waste
$bhrbef.Range("B17").Formula = "='US values'!B17"
$bhrbef.Range("C17").Formula = "='US values'!C17"
$bhrbef.Range("D17").Formula = "='US values'!D17"

# ---------------------------------------------------------------------
# 3) Document the US-values fallback on the "About" sheet (row 42).
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("A42").Value = "Where there is no Korea data available, we use US values from the US EPS."

# Leave the "About" sheet as the active/selected tab (matches the saved
# view state of the edited workbook).
$about.Activate()
